$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-29: columns C-G values updated (norm updates), year in column B unchanged ---
$ws.Range("C2").Value = 214592
$ws.Range("D2").Value = 23037
$ws.Range("E2").Value = 1956
$ws.Range("F2").Value = 11363
$ws.Range("G2").Value = 29829

$ws.Range("C3").Value = 207134
$ws.Range("D3").Value = 19479
$ws.Range("E3").Value = 10002
$ws.Range("F3").Value = 10139
$ws.Range("G3").Value = 34023

$ws.Range("C4").Value = 203183
$ws.Range("D4").Value = 17740
$ws.Range("E4").Value = 14202
$ws.Range("F4").Value = 9521
$ws.Range("G4").Value = 36131

$ws.Range("C5").Value = 198687
$ws.Range("D5").Value = 15971
$ws.Range("E5").Value = 19100
$ws.Range("F5").Value = 8896
$ws.Range("G5").Value = 38123

$ws.Range("C6").Value = 193971
$ws.Range("D6").Value = 14209
$ws.Range("E6").Value = 24221
$ws.Range("F6").Value = 8300
$ws.Range("G6").Value = 40076

$ws.Range("C7").Value = 189006
$ws.Range("D7").Value = 12429
$ws.Range("E7").Value = 29670
$ws.Range("F7").Value = 7712
$ws.Range("G7").Value = 41960

$ws.Range("C8").Value = 183519
$ws.Range("D8").Value = 10710
$ws.Range("E8").Value = 35398
$ws.Range("F8").Value = 7197
$ws.Range("G8").Value = 43953

$ws.Range("C9").Value = 178004
$ws.Range("D9").Value = 8874
$ws.Range("E9").Value = 41352
$ws.Range("F9").Value = 6597
$ws.Range("G9").Value = 45950

$ws.Range("C10").Value = 172286
$ws.Range("D10").Value = 7087
$ws.Range("E10").Value = 47489
$ws.Range("F10").Value = 6030
$ws.Range("G10").Value = 47885

$ws.Range("C11").Value = 166272
$ws.Range("D11").Value = 5299
$ws.Range("E11").Value = 53935
$ws.Range("F11").Value = 5502
$ws.Range("G11").Value = 49769

$ws.Range("C12").Value = 159990
$ws.Range("D12").Value = 3539
$ws.Range("E12").Value = 60644
$ws.Range("F12").Value = 4904
$ws.Range("G12").Value = 51700

$ws.Range("C13").Value = 153524
$ws.Range("D13").Value = 1792
$ws.Range("E13").Value = 67657
$ws.Range("F13").Value = 4344
$ws.Range("G13").Value = 53460

$ws.Range("C14").Value = 146888
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 74835
$ws.Range("F14").Value = 3751
$ws.Range("G14").Value = 55303

$ws.Range("C15").Value = 141364
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 81459
$ws.Range("F15").Value = 3029
$ws.Range("G15").Value = 54925

$ws.Range("C16").Value = 138260
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 85640
$ws.Range("F16").Value = 2326
$ws.Range("G16").Value = 54551

$ws.Range("C17").Value = 135370
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 89700
$ws.Range("F17").Value = 1590
$ws.Range("G17").Value = 54117

$ws.Range("C18").Value = 136929
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 88802
$ws.Range("F18").Value = 1398
$ws.Range("G18").Value = 53648

$ws.Range("C19").Value = 134788
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 91585
$ws.Range("F19").Value = 1306
$ws.Range("G19").Value = 53098

$ws.Range("C20").Value = 132010
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 94921
$ws.Range("F20").Value = 1224
$ws.Range("G20").Value = 52622

$ws.Range("C21").Value = 128991
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 98496
$ws.Range("F21").Value = 1137
$ws.Range("G21").Value = 52153

$ws.Range("C22").Value = 125805
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 102251
$ws.Range("F22").Value = 1040
$ws.Range("G22").Value = 51681

$ws.Range("C23").Value = 122405
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 106258
$ws.Range("F23").Value = 939
$ws.Range("G23").Value = 51175

$ws.Range("C24").Value = 118944
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 110286
$ws.Range("F24").Value = 863
$ws.Range("G24").Value = 50684

$ws.Range("C25").Value = 115423
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 114314
$ws.Range("F25").Value = 778
$ws.Range("G25").Value = 50262

$ws.Range("C26").Value = 111961
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 118311
$ws.Range("F26").Value = 676
$ws.Range("G26").Value = 49829

$ws.Range("C27").Value = 103623
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 127060
$ws.Range("F27").Value = 603
$ws.Range("G27").Value = 49491

$ws.Range("C28").Value = 99939
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 131353
$ws.Range("F28").Value = 513
$ws.Range("G28").Value = 48972

$ws.Range("C29").Value = 96553
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 135311
$ws.Range("F29").Value = 418
$ws.Range("G29").Value = 48495

# --- Rows 30-32: initialisation-order fix -> column B now holds small sequential id,
#     C-F re-derived, and new column G value added (previously absent) ---
$ws.Range("B30").Value = 2066
$ws.Range("C30").Value = 53769
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 170079
$ws.Range("F30").Value = 7415
$ws.Range("G30").Value = 49514

$ws.Range("B31").Value = 2067
$ws.Range("C31").Value = 52816
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 171242
$ws.Range("F31").Value = 7413
$ws.Range("G31").Value = 49306

$ws.Range("B32").Value = 2068
$ws.Range("C32").Value = 51886
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 172402
$ws.Range("F32").Value = 7413
$ws.Range("G32").Value = 49076

# --- Rows 33-41: new simulation years appended, columns B-G only (no column A) ---
$ws.Range("B33").Value = 2069
$ws.Range("C33").Value = 51044
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 173462
$ws.Range("F33").Value = 7389
$ws.Range("G33").Value = 48882

$ws.Range("B34").Value = 2070
$ws.Range("C34").Value = 50250
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 174459
$ws.Range("F34").Value = 7391
$ws.Range("G34").Value = 48677

$ws.Range("B35").Value = 2071
$ws.Range("C35").Value = 49428
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 175512
$ws.Range("F35").Value = 7384
$ws.Range("G35").Value = 48453

$ws.Range("B36").Value = 2072
$ws.Range("C36").Value = 48024
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 177088
$ws.Range("F36").Value = 7350
$ws.Range("G36").Value = 48315

$ws.Range("B37").Value = 2073
$ws.Range("C37").Value = 47326
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 177994
$ws.Range("F37").Value = 7309
$ws.Range("G37").Value = 48148

$ws.Range("B38").Value = 2074
$ws.Range("C38").Value = 46685
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 178858
$ws.Range("F38").Value = 7265
$ws.Range("G38").Value = 47969

$ws.Range("B39").Value = 2075
$ws.Range("C39").Value = 46068
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 179684
$ws.Range("F39").Value = 7222
$ws.Range("G39").Value = 47803

$ws.Range("B40").Value = 2076
$ws.Range("C40").Value = 45498
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 180470
$ws.Range("F40").Value = 7191
$ws.Range("G40").Value = 47618

$ws.Range("B41").Value = 2077
$ws.Range("C41").Value = 44959
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 181232
$ws.Range("F41").Value = 7154
$ws.Range("G41").Value = 47432
